$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dwellings_buildings")

# Row 2 (H:1, Wholesale and retail trade)
$ws.Cells.Item(2, 3).Value = 130
$ws.Cells.Item(2, 4).Value = 287.5895765472313

# Row 3 (H:1, Financial intermediation... -> Professional and technical services)
$ws.Cells.Item(3, 2).Value = "Professional and technical services"
$ws.Cells.Item(3, 3).Value = 130
$ws.Cells.Item(3, 4).Value = 330.3322475570033

# Row 4 (H:1, Community; social and personal services -> All other services)
$ws.Cells.Item(4, 2).Value = "All other services"
$ws.Cells.Item(4, 3).Value = 130
$ws.Cells.Item(4, 4).Value = 298.14332247557

# Row 5 (H:2, Wholesale and retail trade)
$ws.Cells.Item(5, 3).Value = 260
$ws.Cells.Item(5, 4).Value = 287.5895765472313
$ws.Cells.Item(5, 6).Value = 2

# Row 6 (H:2, Financial intermediation... -> Professional and technical services)
$ws.Cells.Item(6, 2).Value = "Professional and technical services"
$ws.Cells.Item(6, 3).Value = 260
$ws.Cells.Item(6, 4).Value = 330.3322475570033
$ws.Cells.Item(6, 6).Value = 2

# Row 7 (H:2, Community; social and personal services -> All other services)
$ws.Cells.Item(7, 2).Value = "All other services"
$ws.Cells.Item(7, 3).Value = 260
$ws.Cells.Item(7, 4).Value = 298.14332247557
$ws.Cells.Item(7, 6).Value = 2

# Row 8 (H:3, Wholesale and retail trade)
$ws.Cells.Item(8, 3).Value = 450
$ws.Cells.Item(8, 4).Value = 287.5895765472313
$ws.Cells.Item(8, 6).Value = 3

# Row 9 (H:3, Financial intermediation... -> Professional and technical services)
$ws.Cells.Item(9, 2).Value = "Professional and technical services"
$ws.Cells.Item(9, 3).Value = 450
$ws.Cells.Item(9, 4).Value = 330.3322475570033
$ws.Cells.Item(9, 6).Value = 3

# Row 10 (H:3, Community; social and personal services -> All other services)
$ws.Cells.Item(10, 2).Value = "All other services"
$ws.Cells.Item(10, 3).Value = 450
$ws.Cells.Item(10, 4).Value = 298.14332247557
$ws.Cells.Item(10, 6).Value = 3

# Row 11 (HBET:3-6, Wholesale and retail trade)
$ws.Cells.Item(11, 3).Value = 900
$ws.Cells.Item(11, 4).Value = 287.5895765472313
$ws.Cells.Item(11, 6).Value = 5

# Row 12 (HBET:3-6, Financial intermediation... -> Professional and technical services)
$ws.Cells.Item(12, 2).Value = "Professional and technical services"
$ws.Cells.Item(12, 3).Value = 900
$ws.Cells.Item(12, 4).Value = 330.3322475570033
$ws.Cells.Item(12, 6).Value = 5

# Row 13 (HBET:3-6, Community; social and personal services -> All other services)
$ws.Cells.Item(13, 2).Value = "All other services"
$ws.Cells.Item(13, 3).Value = 900
$ws.Cells.Item(13, 4).Value = 298.14332247557
$ws.Cells.Item(13, 6).Value = 5

# Row 14 (HBET:4-7, Financial intermediation... -> Professional and technical services)
$ws.Cells.Item(14, 2).Value = "Professional and technical services"
$ws.Cells.Item(14, 3).Value = 1200
$ws.Cells.Item(14, 4).Value = 432.7035830618893
$ws.Cells.Item(14, 6).Value = 5

# Row 15 (HBET:4-7, Wholesale and retail trade -> Professional and technical services)
$ws.Cells.Item(15, 2).Value = "Professional and technical services"
$ws.Cells.Item(15, 3).Value = 1200
$ws.Cells.Item(15, 4).Value = 287.5895765472313
$ws.Cells.Item(15, 6).Value = 5

# Row 16 (HBET:4-7, Community; social and personal services -> All other services)
$ws.Cells.Item(16, 2).Value = "All other services"
$ws.Cells.Item(16, 3).Value = 1200
$ws.Cells.Item(16, 4).Value = 298.14332247557
$ws.Cells.Item(16, 6).Value = 5

# Row 17 (HBET:8+, Financial intermediation... -> Professional and technical services)
$ws.Cells.Item(17, 2).Value = "Professional and technical services"
$ws.Cells.Item(17, 3).Value = 3200
$ws.Cells.Item(17, 4).Value = 432.7035830618893
$ws.Cells.Item(17, 6).Value = 10

# Row 18 (HBET:8+, Wholesale and retail trade -> Professional and technical services)
$ws.Cells.Item(18, 2).Value = "Professional and technical services"
$ws.Cells.Item(18, 3).Value = 3200
$ws.Cells.Item(18, 4).Value = 287.5895765472313
$ws.Cells.Item(18, 6).Value = 10

# Row 19 (HBET:8+, Community; social and personal services -> All other services)
$ws.Cells.Item(19, 2).Value = "All other services"
$ws.Cells.Item(19, 3).Value = 3200
$ws.Cells.Item(19, 4).Value = 298.14332247557
$ws.Cells.Item(19, 6).Value = 10
